$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "561.67") are written as TEXT,
# matching the source inlineStr cells, not auto-converted to numbers.

$ws.Range("D2").Value = '58.987.53'
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = '2.987.11'
$ws.Range("E3").Value = '  +3.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.67'
$ws.Range("E5").Value = '  +1.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.42'
$ws.Range("E6").Value = '  +11.60%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  +4.54%  '

$ws.Range("D9").Value = '2.975.05'
$ws.Range("E9").Value = '  +2.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.132'
$ws.Range("E10").Value = '  +6.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.87'
$ws.Range("E11").Value = '  +3.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").Value = '  +3.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +7.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.59'
$ws.Range("E14").Value = '  +3.42%  '

$ws.Range("E15").Value = '  +2.93%  '

$ws.Range("D16").Value = '3.477.72'
$ws.Range("E16").Value = '  +3.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.98'
$ws.Range("E17").Value = '  +6.58%  '

$ws.Range("D18").Value = '2.976.99'
$ws.Range("E18").Value = '  +2.82%  '

$ws.Range("D19").Value = '58.830.36'
$ws.Range("E19").Value = '  +2.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '423.22'
$ws.Range("E20").Value = '  +4.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.53'
$ws.Range("E21").Value = '  +4.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("E22").Value = '  +5.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.13'
$ws.Range("E23").Value = '  +4.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.45'
$ws.Range("E24").Value = '  +5.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.33'
$ws.Range("E25").Value = '  +4.27%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.12'
$ws.Range("E28").Value = '  +9.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.53'
$ws.Range("E29").Value = '  +3.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.76'
$ws.Range("E30").Value = '  +7.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.59'
$ws.Range("E31").Value = '  +3.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.10'
$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0991'
$ws.Range("E33").Value = '  +0.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").Value = '  +10.94%  '

$ws.Range("D35").Value = '0.0₃0770'
$ws.Range("E35").Value = '  +23.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("E36").Value = '  +5.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.08'
$ws.Range("E37").Value = '  +3.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.88'
$ws.Range("E38").Value = '  +1.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.65'
$ws.Range("E39").Value = '  +3.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.77'
$ws.Range("E40").Value = '  +14.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '400.73'
$ws.Range("E41").Value = '  +11.15%  '

$ws.Range("D42").Value = '2.744.31'
$ws.Range("E42").Value = '  +4.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0349'
$ws.Range("E43").Value = '  +2.18%  '

$ws.Range("E44").Value = '  +0.67%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.245'
$ws.Range("E45").Value = '  +6.90%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.06'
$ws.Range("E47").Value = '  +5.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.02'
$ws.Range("E48").Value = '  +3.39%  '

$ws.Range("E49").Value = '  +2.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.34'
$ws.Range("E50").Value = '  +20.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.31'
$ws.Range("E51").Value = '  +1.66%  '
